$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold, border, centered) from the H1 header cell
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for columns I and J, rows 2-26
$values = @(
    @(1, 5),
    @(1, 5),
    @(7, 9),
    @(6, 8),
    @(5, 7),
    @(5, 9),
    @(1, 6),
    @(1, 6),
    @(1, 6),
    @(1, 3),
    @(1, 4),
    @(1, 5),
    @(1, 3),
    @(1, 6),
    @(1, 5),
    @(1, 6),
    @(1, 3),
    @(1, 4),
    @(1, 6),
    @(1, 7),
    @(1, 6),
    @(1, 7),
    @(1, 3),
    @(1, 7),
    @(1, 6)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
